$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 297.45587
$ws.Range("J17").Value = 297.45587
$ws.Range("L17").Value = 892.36761
$ws.Range("N17").Value = -1228.36761
$ws.Range("H28").Value = 430.05264
$ws.Range("I28").Value = 481.46155
$ws.Range("J28").Value = 318.66666
$ws.Range("K28").Value = 481.46155
$ws.Range("L28").Value = 318.66666
$ws.Range("M28").Value = 3.538450000000012
$ws.Range("N28").Value = -1288.66666
$ws.Range("H99").Value = 946.9545000000001
$ws.Range("I99").Value = 546.8
$ws.Range("J99").Value = 1804.4286
$ws.Range("K99").Value = 1640.4
$ws.Range("L99").Value = 5413.2858
$ws.Range("M99").Value = -142.3999999999999
$ws.Range("N99").Value = -8409.2858
$ws.Range("H112").Value = 1322.3793
$ws.Range("I112").Value = 618
$ws.Range("J112").Value = 1469.125
$ws.Range("K112").Value = 1854
$ws.Range("L112").Value = 4407.375
$ws.Range("M112").Value = -746
$ws.Range("N112").Value = -6623.375
$ws.Range("H125").Value = 2650.5881
$ws.Range("I125").Value = 2020.6666
$ws.Range("J125").Value = 2785.5715
$ws.Range("K125").Value = 18185.9994
$ws.Range("L125").Value = 25070.1435
$ws.Range("M125").Value = -15725.9994
$ws.Range("N125").Value = -29990.1435
$ws.Range("H127").Value = 1021.2
$ws.Range("I127").Value = 600
$ws.Range("J127").Value = 1201.7142
$ws.Range("K127").Value = 1800
$ws.Range("L127").Value = 3605.1426
$ws.Range("M127").Value = 3160
$ws.Range("N127").Value = -13525.1426
$ws.Range("H138").Value = 2884.0815
$ws.Range("I138").Value = 2135
$ws.Range("J138").Value = 3445.8928
$ws.Range("K138").Value = 6405
$ws.Range("L138").Value = 10337.6784
$ws.Range("M138").Value = -1265
$ws.Range("N138").Value = -20617.6784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1226.5
$ws.Range("I45").Value = 958.63635
$ws.Range("J45").Value = 1494.3636
$ws.Range("K45").Value = 958.63635
$ws.Range("L45").Value = 1494.3636
$ws.Range("M45").Value = -581.63635
$ws.Range("N45").Value = -2248.3636
$ws.Range("H122").Value = 2818.9429
$ws.Range("I122").Value = 2271.5806
$ws.Range("J122").Value = 7061
$ws.Range("K122").Value = 6814.7418
$ws.Range("L122").Value = 21183
$ws.Range("M122").Value = -4364.7418
$ws.Range("N122").Value = -26083

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4349730
$ws.Range("I105").Value = 1927
$ws.Range("J105").Value = 7694194
$ws.Range("K105").Value = 1927
$ws.Range("L105").Value = 7694194
$ws.Range("M105").Value = -180
$ws.Range("N105").Value = -7697688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 664.1667
$ws.Range("I5").Value = 92.333336
$ws.Range("J5").Value = 1236
$ws.Range("K5").Value = 92.333336
$ws.Range("L5").Value = 1236
$ws.Range("M5").Value = 19.666664
$ws.Range("N5").Value = -1460
$ws.Range("H31").Value = 2285.9033
$ws.Range("I31").Value = 1596.0238
$ws.Range("J31").Value = 3734.65
$ws.Range("K31").Value = 1596.0238
$ws.Range("L31").Value = 3734.65
$ws.Range("M31").Value = -1301.0238
$ws.Range("N31").Value = -4324.65
$ws.Range("H34").Value = 2285.9033
$ws.Range("I34").Value = 1596.0238
$ws.Range("J34").Value = 3734.65
$ws.Range("K34").Value = 1596.0238
$ws.Range("L34").Value = 3734.65
$ws.Range("M34").Value = -1394.0238
$ws.Range("N34").Value = -4138.65
$ws.Range("H111").Value = 28000
$ws.Range("J111").Value = 28000
$ws.Range("L111").Value = 28000
$ws.Range("N111").Value = -36180
$ws.Range("H134").Value = 1414.091
$ws.Range("I134").Value = 934.6177
$ws.Range("J134").Value = 3044.3
$ws.Range("K134").Value = 2803.8531
$ws.Range("L134").Value = 9132.900000000001
$ws.Range("M134").Value = -268.8531000000003
$ws.Range("N134").Value = -14202.9
$ws.Range("H141").Value = 34496.25
$ws.Range("J141").Value = 34496.25
$ws.Range("L141").Value = 34496.25
$ws.Range("N141").Value = -44856.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 8.307693
$ws.Range("I12").Value = 18.88889
$ws.Range("J12").Value = 2.7058823
$ws.Range("K12").Value = 56.66667
$ws.Range("L12").Value = 8.1176469
$ws.Range("M12").Value = 116.33333
$ws.Range("N12").Value = -354.1176469
$ws.Range("H104").Value = 2153.2
$ws.Range("J104").Value = 2153.2
$ws.Range("L104").Value = 6459.599999999999
$ws.Range("N104").Value = -11701.6
$ws.Range("H121").Value = 4000.3333
$ws.Range("I121").Value = 6822
$ws.Range("J121").Value = 3118.5625
$ws.Range("K121").Value = 20466
$ws.Range("L121").Value = 9355.6875
$ws.Range("M121").Value = -19156
$ws.Range("N121").Value = -11975.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4744.381
$ws.Range("I102").Value = 3306
$ws.Range("K102").Value = 3306
$ws.Range("M102").Value = -1684
$ws.Range("H122").Value = 1179.8125
$ws.Range("I122").Value = 1172.4445
$ws.Range("J122").Value = 1189.2858
$ws.Range("K122").Value = 3517.3335
$ws.Range("L122").Value = 3567.8574
$ws.Range("M122").Value = -1067.3335
$ws.Range("N122").Value = -8467.857400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 826
$ws.Range("I61").Value = 826
$ws.Range("K61").Value = 826
$ws.Range("M61").Value = -624
$ws.Range("H82").Value = 2342.2727
$ws.Range("I82").Value = 1756
$ws.Range("J82").Value = 3368.25
$ws.Range("K82").Value = 1756
$ws.Range("L82").Value = 3368.25
$ws.Range("M82").Value = -1395
$ws.Range("N82").Value = -4090.25
$ws.Range("H85").Value = 2342.2727
$ws.Range("I85").Value = 1756
$ws.Range("J85").Value = 3368.25
$ws.Range("K85").Value = 1756
$ws.Range("L85").Value = 3368.25
$ws.Range("M85").Value = -508
$ws.Range("N85").Value = -5864.25
$ws.Range("H100").Value = 47622460
$ws.Range("I100").Value = 3981.8823
$ws.Range("J100").Value = 250001000
$ws.Range("K100").Value = 3981.8823
$ws.Range("L100").Value = 250001000
$ws.Range("M100").Value = -3440.8823
$ws.Range("N100").Value = -250002082
$ws.Range("H113").Value = 826
$ws.Range("I113").Value = 826
$ws.Range("K113").Value = 826
$ws.Range("M113").Value = 1344
$ws.Range("H122").Value = 2389.3635
$ws.Range("I122").Value = 2531.2778
$ws.Range("J122").Value = 1750.75
$ws.Range("K122").Value = 7593.8334
$ws.Range("L122").Value = 5252.25
$ws.Range("M122").Value = -5143.8334
$ws.Range("N122").Value = -10152.25
$ws.Range("H132").Value = 7777.8647
$ws.Range("I132").Value = 2152.762
$ws.Range("J132").Value = 15160.8125
$ws.Range("K132").Value = 6458.286
$ws.Range("L132").Value = 45482.4375
$ws.Range("M132").Value = -3928.286
$ws.Range("N132").Value = -50542.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1038.08
$ws.Range("I126").Value = 377.9091
$ws.Range("J126").Value = 1556.7858
$ws.Range("K126").Value = 1133.7273
$ws.Range("L126").Value = 4670.357400000001
$ws.Range("M126").Value = 1336.2727
$ws.Range("N126").Value = -9610.357400000001
